# Rename Gaussian methods to reflect underlying algorithms.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmarks")

# Rename the method/algorithm labels used in both the 2048x2048 benchmark
# block (rows 20-23) and the 512x512 benchmark block (rows 28-31).
$ws.Range("C20").Value = "GaussianFilter"
$ws.Range("C21").Value = "GaussianBlur"
$ws.Range("C22").Value = "GaussianBlurBox"
$ws.Range("C23").Value = "GaussianBlurBoxIndependent"

$ws.Range("C28").Value = "GaussianFilter"
$ws.Range("C29").Value = "GaussianBlur"
$ws.Range("C30").Value = "GaussianBlurBox"
$ws.Range("C31").Value = "GaussianBlurBoxIndependent"

# Updated timing measurements for the 512x512 benchmark.
$ws.Range("I28").Value = 8.9
$ws.Range("I29").Value = 11.017
$ws.Range("I31").Value = 11.3

# Move / extend the active selection to the updated GaussianBlur rows.
[void]$ws.Range("C28:C31").Select()
